$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Range("A1").Value = 45311
$ws.Range("D28").Value = 469.482
$ws.Range("D29").Value = 536.554
$ws.Range("D30").Value = 613.997
$ws.Range("D31").Value = 718.271
$ws.Range("D32").Value = 785.847
$ws.Range("D33").Value = 936.438
$ws.Range("D34").Value = 959.213
$ws.Range("D35").Value = 1070.066
$ws.Range("D36").Value = 1295.822
$ws.Range("D37").Value = 1462.864
$ws.Range("D38").Value = 1708.361
$ws.Range("D39").Value = 1885.526
$ws.Range("D40").Value = 2050.034
$ws.Range("D41").Value = 2303.124
$ws.Range("D42").Value = 2457.509
$ws.Range("D43").Value = 2680.229
$ws.Range("D49").Value = 1505.89
$ws.Range("D50").Value = 1513.482
$ws.Range("D51").Value = 1857.684
$ws.Range("D52").Value = 1986.763
$ws.Range("D53").Value = 2116.086
$ws.Range("D54").Value = 2381.58
$ws.Range("D55").Value = 2771.34
$ws.Range("D56").Value = 3102.889
$ws.Range("D57").Value = 3548.327
$ws.Range("D58").Value = 3998.828
$ws.Range("D59").Value = 4340.497
$ws.Range("D60").Value = 4884.644
$ws.Range("D61").Value = 5388.291
$ws.Range("D62").Value = 5707.184
$ws.Range("D68").Value = 607.671
$ws.Range("D69").Value = 625.134
$ws.Range("D70").Value = 707.388
$ws.Range("D71").Value = 824.064
$ws.Range("D72").Value = 953.141
$ws.Range("D73").Value = 1104.993
$ws.Range("D74").Value = 1237.612
$ws.Range("D75").Value = 1311.008
$ws.Range("D76").Value = 1543.85
$ws.Range("D77").Value = 1789.351
$ws.Range("D78").Value = 1993.592
$ws.Range("D79").Value = 2234.28
$ws.Range("D80").Value = 2492.942
$ws.Range("D81").Value = 2640.49
$ws.Range("D82").Value = 2862.453
$ws.Range("D83").Value = 3125.665
$ws.Range("D89").Value = 259.878
$ws.Range("D90").Value = 312.062
$ws.Range("D91").Value = 337.37
$ws.Range("D92").Value = 418.61
$ws.Range("D93").Value = 456.827
$ws.Range("D94").Value = 492.008
$ws.Range("D95").Value = 539.842
$ws.Range("D96").Value = 688.913
$ws.Range("D97").Value = 749.146
$ws.Range("D98").Value = 878.223
$ws.Range("D99").Value = 992.114
$ws.Range("D100").Value = 1068.042
$ws.Range("D106").Value = 2406.888
$ws.Range("D107").Value = 2558.744
$ws.Range("D108").Value = 2708.068
$ws.Range("D109").Value = 2961.16
$ws.Range("D110").Value = 2978.875
$ws.Range("D111").Value = 3510.364
$ws.Range("D112").Value = 4155.742
$ws.Range("D113").Value = 4606.245
$ws.Range("D114").Value = 5036.501
$ws.Range("D115").Value = 5560.397
$ws.Range("D116").Value = 6226.021
$ws.Range("D117").Value = 6577.82
$ws.Range("D118").Value = 7339.619
$ws.Range("D119").Value = 7782.528
